$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Runmode column (C) for rows 3-5 from "N" to "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"

# Update the active selection to E17
$ws.Range("E17").Select()
